# poprawka swingu, dodanie wydmuchu
$wb = $excel.ActiveWorkbook

# Duplicate "Zestaw2" to get a sheet with the identical layout/formulas/styles,
# placed right after it (becomes the last, newly-selected tab).
$ws2 = $wb.Worksheets.Item("Zestaw2")
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item("Zestaw2 (2)")
$ws3.Name = "Zestaw3"

# New measurement data for Zestaw3 (column A = KS4 temperature, text-typed;
# columns B..E = numeric readings).
$colA = @("376.7","376.7","376.7","376.7","375","371.5","371.5","371.5","292.5","292.5","292.5","292.5","292.5","295.7","295.7","302.6","297.4","296.1","296.1","296.1","296.1","296.1","296.1","296.1")
$colB = @(45,45,45,45,46,46,45,46,45,45,45,45,45,45,45,46,45,45,45,45,45,45,45,45)
$colC = @(99,99,99,99,99,97,100,102,86,86,86,86,86,86,89,92,87,87,87,87,87,87,90,92)
$colD = @(541,551,551,549,547,530,531,536,509,509,509,509,508,511,512,518,512,508,508,508,508,508,508,517)
$colE = @(141,144,144,144,143,139,139,139,129,129,129,129,129,130,130,132,131,130,130,128,130,130,130,134)

$r = 2
foreach ($val in $colA) {
    # Store as a formula first so it evaluates to a literal text string;
    # converted to a plain value below via copy / paste-special.
    $ws3.Cells.Item($r, 1).Formula = '="' + $val + '"'
    $r++
}
$ws3.Range("A2:A25").Copy()
$ws3.Range("A2:A25").PasteSpecial(-4163)

$r = 2
foreach ($val in $colB) { $ws3.Cells.Item($r, 2).Value = $val; $r++ }
$r = 2
foreach ($val in $colC) { $ws3.Cells.Item($r, 3).Value = $val; $r++ }
$r = 2
foreach ($val in $colD) { $ws3.Cells.Item($r, 4).Value = $val; $r++ }
$r = 2
foreach ($val in $colE) { $ws3.Cells.Item($r, 5).Value = $val; $r++ }

# Row 6, column A is a genuine number (375), not text like the rest of the column.
$ws3.Range("A6").Value = 375

# Selection on the new sheet matches the source workbook (J18); activating it
# also makes it the selected / active tab.
$ws3.Range("J18").Select()
$ws3.Activate()
